$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Testing")

# Rows whose Test Status / Requirement Status become starred variants ("Pass*"/"Complete*").
# These correspond to TestID 202, 204, 205, 206, 207, 209 (rows 14, 16, 17, 18, 19, 21).
$starredRows = @(14, 16, 17, 18, 19, 21)
foreach ($r in $starredRows) {
    $ws2.Cells.Item($r, 3).Value = "Pass*"
    $ws2.Cells.Item($r, 7).Value = "Complete*"
}

# Update the workbook tab scroll position (firstSheet) - best effort.
$win = $wb.Windows.Item(1)
$win.ScrollWorkbookTabs(1, 1)

# Activate the Testing sheet and move the selection/active cell to E18, matching the new view state.
$ws2.Activate()
$ws2.Range("E18").Select()
